# C5-PowerPoint.pptx edit: table style id swap + theme colour scheme swap
# (theme1.xml "Integral" colours -> "Office" colours; the notes-master
# theme part is not independently addressable through this COM host, so
# only the slide-master theme (ppt/theme/theme1.xml) is updated here.)

$p = $ppt.ActivePresentation

# --- 1. Table on slide 6: swap the applied table style GUID -------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{5EDC193C-A99A-43F3-9082-81FA7BC53D57}")
}

# --- 2. Slide master theme: recolour from "Integral" to "Office" --------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

# index : slot      : target RGB (hex)
#   1   : dk1        : 000000
#   2   : lt1        : FFFFFF
#   3   : dk2        : 44546A
#   4   : lt2        : E7E6E6
#   5   : accent1    : 5B9BD5
#   6   : accent2    : ED7D31
#   7   : accent3    : A5A5A5
#   8   : accent4    : FFC000
#   9   : accent5    : 4472C4
#  10   : accent6    : 70AD47
#  11   : hlink      : 0563C1
#  12   : folHlink   : 954F72
$officeRgb = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeRgb[$i - 1]
}
